# Regenerate orders with updated distance/size codes.
# Distance codes: D64 -> D69, D80 -> D86, D51 -> D55
# Size code:      S30 -> S31
# These substitutions are applied across every cell in the sheet (Trial
# labels, Filename_Left / Filename_Right, Distance, Size columns, etc.)
# so every string built from the old codes picks up the new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
